$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I (I0) and J (IF), styled like the existing H1 header
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2..34
$data = @{
    2  = @(1, 7)
    3  = @(1, 7)
    4  = @(1, 7)
    5  = @(1, 6)
    6  = @(1, 6)
    7  = @(7, 8)
    8  = @(1, 5)
    9  = @(1, 7)
    10 = @(1, 2)
    11 = @(1, 5)
    12 = @(1, 6)
    13 = @(1, 7)
    14 = @(1, 6)
    15 = @(4, 9)
    16 = @(1, 5)
    17 = @(1, 4)
    18 = @(2, 2)
    19 = @(4, 6)
    20 = @(3, 5)
    21 = @(1, 5)
    22 = @(1, 4)
    23 = @(1, 5)
    24 = @(1, 6)
    25 = @(1, 5)
    26 = @(1, 6)
    27 = @(4, 8)
    28 = @(2, 7)
    29 = @(9, 9)
    30 = @(5, 6)
    31 = @(1, 3)
    32 = @(4, 5)
    33 = @(3, 4)
    34 = @(3, 4)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 9).Value = $values[0]
    $ws.Cells.Item($row, 10).Value = $values[1]
}
